$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.321.95'
$ws.Range('E2').Value = '  -0.16%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.430.09'
$ws.Range('E3').Value = '  +1.13%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.24'
$ws.Range('E5').Value = '  -1.44%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.29'
$ws.Range('E6').Value = '  -1.99%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.422.98'
$ws.Range('E7').Value = '  +1.13%  '

$ws.Range('E8').Value = '  +0.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.591'
$ws.Range('E9').Value = '  -0.80%  '

$ws.Range('E10').Value = '  +0.56%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.585'
$ws.Range('E11').Value = '  -0.76%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '48.74'
$ws.Range('E12').Value = '  -0.32%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000280'
$ws.Range('E13').Value = '  -1.61%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '699.23'
$ws.Range('E14').Value = '  +2.03%  '

$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.964.17'
$ws.Range('E15').Value = '  +0.60%  '

$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.63'
$ws.Range('E16').Value = '  +0.25%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.429.94'
$ws.Range('E17').Value = '  -0.05%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.431.83'
$ws.Range('E18').Value = '  +1.47%  '

$ws.Range('E19').Value = '  +0.78%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.71'
$ws.Range('E20').Value = '  -0.08%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.41'
$ws.Range('E21').Value = '  -0.04%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.898'
$ws.Range('E22').Value = '  -0.72%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.39'
$ws.Range('E23').Value = '  -0.11%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '16.92'
$ws.Range('E24').Value = '  -1.09%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '100.90'
$ws.Range('E25').Value = '  -3.56%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.90'
$ws.Range('E26').Value = '  -1.68%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.67'
$ws.Range('E27').Value = '  -2.13%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.60'
$ws.Range('E28').Value = '  -0.41%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.52'
$ws.Range('E29').Value = '  -3.31%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.74'
$ws.Range('E30').Value = '  +0.29%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.97'
$ws.Range('E31').Value = '  -1.26%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.73'
$ws.Range('E32').Value = '  +1.57%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '567.60'
$ws.Range('E33').Value = '  +2.16%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.02'
$ws.Range('E34').Value = '  -1.53%  '

$ws.Range('E35').Value = '  -1.59%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.22'
$ws.Range('E36').Value = '  -0.33%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.04%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.601.58'
$ws.Range('E38').Value = '  -3.28%  '

$ws.Range('E39').Value = '  -1.92%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '34.94'
$ws.Range('E40').Value = '  -0.29%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₃0728'
$ws.Range('E41').Value = '  +2.09%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.28'
$ws.Range('E42').Value = '  +0.83%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.67'
$ws.Range('E43').Value = '  -0.18%  '

$ws.Range('E44').Value = '  +3.17%  '

$ws.Range('E45').Value = '  -2.26%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0419'
$ws.Range('E46').Value = '  -0.04%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.48'
$ws.Range('E47').Value = '  +4.41%  '

$ws.Range('E48').Value = '  -0.19%  '

$ws.Range('E49').Value = '  -1.46%  '

$ws.Range('E50').Value = '  -0.17%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '131.44'
$ws.Range('E51').Value = '  -1.00%  '
